# Add 5 new running examples (rows 5-9) to the Laufdaten worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: date (serial), distance (B), hours (C), minutes (D), seconds (E)
$rows = @(
    @{ Row=5; Date="2022-07-24"; B=3;       C=0; D=11; E=45 },
    @{ Row=6; Date="2022-07-25"; B=2;       C=0; D=6;  E=15 },
    @{ Row=7; Date="2022-07-27"; B=5;       C=0; D=30; E=0  },
    @{ Row=8; Date="2022-07-28"; B=10;      C=0; D=55; E=0  },
    @{ Row=9; Date="2022-07-30"; B=21.0975; C=1; D=35; E=0  }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E

    $ws.Cells.Item($row, 7).Formula = "=(C$row*60)+D$row+(E$row/60)"
    $ws.Cells.Item($row, 8).Formula = "=G$row/B$row"
    $ws.Cells.Item($row, 9).Formula = "=B$row/(G$row/60)"
}

$ws.Range("H2").Select()
